$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $value) {
    $c = $ws.Range($addr)
    $origStyle = $c.Style
    $c.NumberFormat = "@"
    $c.Value = $value
    $c.Style = $origStyle
}

Set-TextValue "D2" "26.967.05"
Set-TextValue "E2" "  +1.90%  "
Set-TextValue "D3" "1.840.71"
Set-TextValue "E3" "  +1.60%  "
Set-TextValue "D4" "1.009"
Set-TextValue "E4" "  +0.39%  "
Set-TextValue "D5" "309.54"
Set-TextValue "E5" "  +1.16%  "
Set-TextValue "E6" "  +0.26%  "
Set-TextValue "D7" "0.4667"
Set-TextValue "E7" "  +3.67%  "
Set-TextValue "D8" "0.3619"
Set-TextValue "E8" "  +1.01%  "
Set-TextValue "D9" "0.07129"
Set-TextValue "E9" "  +1.01%  "
Set-TextValue "D10" "0.9111"
Set-TextValue "E10" "  +2.22%  "
Set-TextValue "D11" "19.52"
Set-TextValue "E11" "  +0.78%  "
Set-TextValue "D12" "0.07687"
Set-TextValue "E12" "  -1.41%  "
Set-TextValue "D13" "1.828.15"
Set-TextValue "E13" "  +0.47%  "
Set-TextValue "D14" "5.269"
Set-TextValue "E14" "  -0.10%  "
Set-TextValue "D15" "6.393"
Set-TextValue "E15" "  +1.25%  "
Set-TextValue "D16" "88.22"
Set-TextValue "E16" "  +3.99%  "
Set-TextValue "E17" "  +0.28%  "
Set-TextValue "D18" "0.000008581"
Set-TextValue "E18" "  +0.59%  "
Set-TextValue "D19" "1.007"
Set-TextValue "D20" "27.003.02"
Set-TextValue "E20" "  +1.97%  "
Set-TextValue "D21" "14.30"
Set-TextValue "E21" "  +0.81%  "
Set-TextValue "D22" "5.012"
Set-TextValue "E22" "  +0.95%  "
Set-TextValue "D23" "10.63"
Set-TextValue "E23" "  +1.05%  "
Set-TextValue "D24" "1.931"
Set-TextValue "E24" "  -1.16%  "
Set-TextValue "D25" "152.44"
Set-TextValue "E25" "  +0.62%  "
Set-TextValue "D26" "18.21"
Set-TextValue "E26" "  +2.41%  "
Set-TextValue "D27" "2.032"
Set-TextValue "E27" "  -1.48%  "
Set-TextValue "D28" "113.99"
Set-TextValue "E28" "  +1.59%  "
Set-TextValue "D29" "4.887"
Set-TextValue "E29" "  +0.67%  "
Set-TextValue "D30" "0.08853"
Set-TextValue "E30" "  +1.94%  "
Set-TextValue "D31" "3.199"
Set-TextValue "E31" "  +2.73%  "
Set-TextValue "D32" "2.835"
Set-TextValue "E32" "  +2.91%  "
Set-TextValue "D33" "0.7481"
Set-TextValue "E33" "  +1.21%  "
Set-TextValue "D34" "1.170"
Set-TextValue "E34" "  +5.24%  "
Set-TextValue "D35" "4.459"
Set-TextValue "E35" "  +0.03%  "
Set-TextValue "D37" "2.983"
Set-TextValue "E37" "  +3.17%  "
Set-TextValue "D38" "0.01937"
Set-TextValue "E38" "  +0.65%  "
Set-TextValue "D39" "0.05161"
Set-TextValue "E39" "  +0.68%  "
Set-TextValue "D40" "0.5175"
Set-TextValue "E40" "  +1.66%  "
Set-TextValue "D41" "6.900"
Set-TextValue "E41" "  +1.92%  "
Set-TextValue "D42" "0.1511"
Set-TextValue "E42" "  +0.29%  "
Set-TextValue "D43" "8.116"
Set-TextValue "E43" "  +0.87%  "
Set-TextValue "D44" "10.46"
Set-TextValue "E44" "  +4.87%  "
Set-TextValue "D45" "0.4692"
Set-TextValue "E45" "  +0.30%  "
Set-TextValue "D46" "1.008"
Set-TextValue "E46" "  +0.20%  "
Set-TextValue "D47" "100.66"
Set-TextValue "E47" "  +0.89%  "
Set-TextValue "D48" "1.603"
Set-TextValue "E48" "  +2.07%  "
Set-TextValue "D49" "0.06043"
Set-TextValue "E49" "  +0.80%  "
Set-TextValue "D50" "64.44"
Set-TextValue "E50" "  +1.38%  "
Set-TextValue "D51" "36.20"
Set-TextValue "E51" "  +0.82%  "
